$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 981.96875
$ws.Range("J17").Value = 1071.2222
$ws.Range("L17").Value = 3213.6666
$ws.Range("N17").Value = -3549.6666
$ws.Range("H32").Value = 586.65
$ws.Range("I32").Value = 424.44446
$ws.Range("J32").Value = 719.36365
$ws.Range("K32").Value = 424.44446
$ws.Range("L32").Value = 719.36365
$ws.Range("M32").Value = -98.44445999999999
$ws.Range("N32").Value = -1371.36365
$ws.Range("H64").Value = 3226.2964
$ws.Range("I64").Value = 3210.4
$ws.Range("J64").Value = 3235.647
$ws.Range("K64").Value = 3210.4
$ws.Range("L64").Value = 3235.647
$ws.Range("M64").Value = -2962.4
$ws.Range("N64").Value = -3731.647
$ws.Range("H67").Value = 3226.2964
$ws.Range("I67").Value = 3210.4
$ws.Range("J67").Value = 3235.647
$ws.Range("K67").Value = 3210.4
$ws.Range("L67").Value = 3235.647
$ws.Range("M67").Value = -2352.4
$ws.Range("N67").Value = -4951.647
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H76").Value = 3654.2856
$ws.Range("I76").Value = 3140
$ws.Range("J76").Value = 4940
$ws.Range("K76").Value = 3140
$ws.Range("L76").Value = 4940
$ws.Range("M76").Value = -2825
$ws.Range("N76").Value = -5570
$ws.Range("H79").Value = 3654.2856
$ws.Range("I79").Value = 3140
$ws.Range("J79").Value = 4940
$ws.Range("K79").Value = 3140
$ws.Range("L79").Value = 4940
$ws.Range("M79").Value = -2048
$ws.Range("N79").Value = -7124
$ws.Range("H138").Value = 2808.9473
$ws.Range("I138").Value = 2374.7646
$ws.Range("J138").Value = 2934.0508
$ws.Range("K138").Value = 7124.293799999999
$ws.Range("L138").Value = 8802.152399999999
$ws.Range("M138").Value = -1984.293799999999
$ws.Range("N138").Value = -19082.1524

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 12771.857
$ws.Range("J95").Value = 12771.857
$ws.Range("L95").Value = 12771.857
$ws.Range("N95").Value = -18263.857
$ws.Range("H101").Value = 13529
$ws.Range("J101").Value = 13529
$ws.Range("L101").Value = 13529
$ws.Range("N101").Value = -20019
$ws.Range("H102").Value = 2725
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 3933.3333
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 3933.3333
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -7177.3333
$ws.Range("H104").Value = 47900.555
$ws.Range("J104").Value = 47900.555
$ws.Range("L104").Value = 47900.555
$ws.Range("N104").Value = -54888.555
$ws.Range("H122").Value = 2260
$ws.Range("I122").Value = 2446.6667
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 7340.000100000001
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -4890.000100000001
$ws.Range("N122").Value = -10000

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2597.4614
$ws.Range("I105").Value = 2436
$ws.Range("J105").Value = 2698.375
$ws.Range("K105").Value = 2436
$ws.Range("L105").Value = 2698.375
$ws.Range("M105").Value = -689
$ws.Range("N105").Value = -6192.375
$ws.Range("H109").Value = 31145
$ws.Range("J109").Value = 31145
$ws.Range("L109").Value = 31145
$ws.Range("N109").Value = -33919
$ws.Range("H134").Value = 2416.5
$ws.Range("I134").Value = 1622.5625
$ws.Range("J134").Value = 3475.0833
$ws.Range("K134").Value = 4867.6875
$ws.Range("L134").Value = 10425.2499
$ws.Range("M134").Value = -2332.6875
$ws.Range("N134").Value = -15495.2499

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 808.6667
$ws.Range("I2").Value = 808.6667
$ws.Range("K2").Value = 808.6667
$ws.Range("M2").Value = -695.6667
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H43").Value = 23414.25
$ws.Range("J43").Value = 23414.25
$ws.Range("L43").Value = 23414.25
$ws.Range("N43").Value = -23782.25
$ws.Range("H62").Value = 9251.066000000001
$ws.Range("J62").Value = 12441.1
$ws.Range("L62").Value = 12441.1
$ws.Range("N62").Value = -13689.1
$ws.Range("H65").Value = 9251.066000000001
$ws.Range("J65").Value = 12441.1
$ws.Range("L65").Value = 62205.5
$ws.Range("N65").Value = -68445.5
$ws.Range("H101").Value = 23414.25
$ws.Range("J101").Value = 23414.25
$ws.Range("L101").Value = 23414.25
$ws.Range("N101").Value = -29904.25
$ws.Range("H131").Value = 21726
$ws.Range("J131").Value = 21726
$ws.Range("L131").Value = 21726
$ws.Range("N131").Value = -31806
$ws.Range("H134").Value = 1474.3096
$ws.Range("I134").Value = 1028.5172
$ws.Range("J134").Value = 2468.7693
$ws.Range("K134").Value = 3085.5516
$ws.Range("L134").Value = 7406.3079
$ws.Range("M134").Value = -550.5515999999998
$ws.Range("N134").Value = -12476.3079

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 5238.615
$ws.Range("I43").Value = 1551
$ws.Range("J43").Value = 5909.091
$ws.Range("K43").Value = 4653
$ws.Range("L43").Value = 17727.273
$ws.Range("M43").Value = -4539
$ws.Range("N43").Value = -17955.273

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7459.048
$ws.Range("I70").Value = 8390
$ws.Range("J70").Value = 6612.727
$ws.Range("K70").Value = 8390
$ws.Range("L70").Value = 6612.727
$ws.Range("M70").Value = -8120
$ws.Range("N70").Value = -7152.727
$ws.Range("H73").Value = 7459.048
$ws.Range("I73").Value = 8390
$ws.Range("J73").Value = 6612.727
$ws.Range("K73").Value = 8390
$ws.Range("L73").Value = 6612.727
$ws.Range("M73").Value = -7454
$ws.Range("N73").Value = -8484.726999999999
$ws.Range("H104").Value = 31313.666
$ws.Range("J104").Value = 31313.666
$ws.Range("L104").Value = 31313.666
$ws.Range("N104").Value = -38301.666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3479.95
$ws.Range("I7").Value = 3709.9
$ws.Range("J7").Value = 3250
$ws.Range("K7").Value = 3709.9
$ws.Range("L7").Value = 3250
$ws.Range("M7").Value = -3597.9
$ws.Range("N7").Value = -3474
$ws.Range("H40").Value = 34485932
$ws.Range("I40").Value = 66669350
$ws.Range("J40").Value = 3693.7856
$ws.Range("K40").Value = 66669350
$ws.Range("L40").Value = 3693.7856
$ws.Range("M40").Value = -66669214
$ws.Range("N40").Value = -3965.7856
$ws.Range("H126").Value = 3479.95
$ws.Range("I126").Value = 3709.9
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 11129.7
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -8659.700000000001
$ws.Range("N126").Value = -14690

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
